$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, "A").Value = "ECs"
$ws.Cells.Item(2, "B").Value = "Egf"
$ws.Cells.Item(2, "C").Value = "Egfr"
$ws.Cells.Item(2, "D").Value = "ECs"
$ws.Cells.Item(2, "E").Value = 2
$ws.Cells.Item(2, "F").Value = 0.6666666666666666
$ws.Cells.Item(2, "G").Value = 0.06368266666666667
$ws.Cells.Item(2, "H").Value = 0.191048
$ws.Cells.Item(2, "I").Value = 0.07976548992723123
$ws.Cells.Item(2, "J").Value = 0.07976548992723123
$ws.Cells.Item(2, "K").Value = 3
$ws.Cells.Item(2, "L").Value = 1
$ws.Cells.Item(2, "M").Value = 1.307106666666667
$ws.Cells.Item(2, "N").Value = 3.92132
$ws.Cells.Item(2, "O").Value = 0.01256263154946851
$ws.Cells.Item(2, "P").Value = 0.01256263154946851
$ws.Cells.Item(2, "Q").Value = 0.08324003815111111
$ws.Cells.Item(2, "R").Value = 0.74916034336
$ws.Cells.Item(2, "S").Value = 0.001002064460318648
$ws.Cells.Item(2, "T").Value = 0.001002064460318648

# Row 3
$ws.Cells.Item(3, "A").Value = "ECs"
$ws.Cells.Item(3, "B").Value = "Egf"
$ws.Cells.Item(3, "C").Value = "Egfr"
$ws.Cells.Item(3, "D").Value = "FAPs"
$ws.Cells.Item(3, "E").Value = 2
$ws.Cells.Item(3, "F").Value = 0.6666666666666666
$ws.Cells.Item(3, "G").Value = 0.06368266666666667
$ws.Cells.Item(3, "H").Value = 0.191048
$ws.Cells.Item(3, "I").Value = 0.07976548992723123
$ws.Cells.Item(3, "J").Value = 0.07976548992723123
$ws.Cells.Item(3, "K").Value = 3
$ws.Cells.Item(3, "L").Value = 1
$ws.Cells.Item(3, "M").Value = 80.22623699999998
$ws.Cells.Item(3, "N").Value = 240.678711
$ws.Cells.Item(3, "O").Value = 0.77105616682495
$ws.Cells.Item(3, "P").Value = 0.77105616682495
$ws.Cells.Item(3, "Q").Value = 5.109020708791999
$ws.Cells.Item(3, "R").Value = 45.98118637912799
$ws.Cells.Item(3, "S").Value = 0.06150367290820508
$ws.Cells.Item(3, "T").Value = 0.06150367290820508

# Row 4
$ws.Cells.Item(4, "A").Value = "ECs"
$ws.Cells.Item(4, "B").Value = "Egf"
$ws.Cells.Item(4, "C").Value = "Egfr"
$ws.Cells.Item(4, "D").Value = "sCs"
$ws.Cells.Item(4, "E").Value = 2
$ws.Cells.Item(4, "F").Value = 0.6666666666666666
$ws.Cells.Item(4, "G").Value = 0.06368266666666667
$ws.Cells.Item(4, "H").Value = 0.191048
$ws.Cells.Item(4, "I").Value = 0.07976548992723123
$ws.Cells.Item(4, "J").Value = 0.07976548992723123
$ws.Cells.Item(4, "K").Value = 3
$ws.Cells.Item(4, "L").Value = 1
$ws.Cells.Item(4, "M").Value = 22.51385866666667
$ws.Cells.Item(4, "N").Value = 67.541576
$ws.Cells.Item(4, "O").Value = 0.2163812016255815
$ws.Cells.Item(4, "P").Value = 0.2163812016255815
$ws.Cells.Item(4, "Q").Value = 1.433742556849778
$ws.Cells.Item(4, "R").Value = 12.903683011648
$ws.Cells.Item(4, "S").Value = 0.01725975255870751
$ws.Cells.Item(4, "T").Value = 0.01725975255870751

# Row 5
$ws.Cells.Item(5, "A").Value = "FAPs"
$ws.Cells.Item(5, "B").Value = "Egf"
$ws.Cells.Item(5, "C").Value = "Egfr"
$ws.Cells.Item(5, "D").Value = "ECs"
$ws.Cells.Item(5, "E").Value = 3
$ws.Cells.Item(5, "F").Value = 1
$ws.Cells.Item(5, "G").Value = 0.3966103333333333
$ws.Cells.Item(5, "H").Value = 1.189831
$ws.Cells.Item(5, "I").Value = 0.4967728144006086
$ws.Cells.Item(5, "J").Value = 0.4967728144006086
$ws.Cells.Item(5, "K").Value = 3
$ws.Cells.Item(5, "L").Value = 1
$ws.Cells.Item(5, "M").Value = 1.307106666666667
$ws.Cells.Item(5, "N").Value = 3.92132
$ws.Cells.Item(5, "O").Value = 0.01256263154946851
$ws.Cells.Item(5, "P").Value = 0.01256263154946851
$ws.Cells.Item(5, "Q").Value = 0.5184120107688889
$ws.Cells.Item(5, "R").Value = 4.66570809692
$ws.Cells.Item(5, "S").Value = 0.00624077383110735
$ws.Cells.Item(5, "T").Value = 0.00624077383110735

# Row 6
$ws.Cells.Item(6, "A").Value = "FAPs"
$ws.Cells.Item(6, "B").Value = "Egf"
$ws.Cells.Item(6, "C").Value = "Egfr"
$ws.Cells.Item(6, "D").Value = "FAPs"
$ws.Cells.Item(6, "E").Value = 3
$ws.Cells.Item(6, "F").Value = 1
$ws.Cells.Item(6, "G").Value = 0.3966103333333333
$ws.Cells.Item(6, "H").Value = 1.189831
$ws.Cells.Item(6, "I").Value = 0.4967728144006086
$ws.Cells.Item(6, "J").Value = 0.4967728144006086
$ws.Cells.Item(6, "K").Value = 3
$ws.Cells.Item(6, "L").Value = 1
$ws.Cells.Item(6, "M").Value = 80.22623699999998
$ws.Cells.Item(6, "N").Value = 240.678711
$ws.Cells.Item(6, "O").Value = 0.77105616682495
$ws.Cells.Item(6, "P").Value = 0.77105616682495
$ws.Cells.Item(6, "Q").Value = 31.81855459864899
$ws.Cells.Item(6, "R").Value = 286.366991387841
$ws.Cells.Item(6, "S").Value = 0.3830397420545756
$ws.Cells.Item(6, "T").Value = 0.3830397420545756

# Row 7
$ws.Cells.Item(7, "A").Value = "FAPs"
$ws.Cells.Item(7, "B").Value = "Egf"
$ws.Cells.Item(7, "C").Value = "Egfr"
$ws.Cells.Item(7, "D").Value = "sCs"
$ws.Cells.Item(7, "E").Value = 3
$ws.Cells.Item(7, "F").Value = 1
$ws.Cells.Item(7, "G").Value = 0.3966103333333333
$ws.Cells.Item(7, "H").Value = 1.189831
$ws.Cells.Item(7, "I").Value = 0.4967728144006086
$ws.Cells.Item(7, "J").Value = 0.4967728144006086
$ws.Cells.Item(7, "K").Value = 3
$ws.Cells.Item(7, "L").Value = 1
$ws.Cells.Item(7, "M").Value = 22.51385866666667
$ws.Cells.Item(7, "N").Value = 67.541576
$ws.Cells.Item(7, "O").Value = 0.2163812016255815
$ws.Cells.Item(7, "P").Value = 0.2163812016255815
$ws.Cells.Item(7, "Q").Value = 8.929228990406223
$ws.Cells.Item(7, "R").Value = 80.36306091365601
$ws.Cells.Item(7, "S").Value = 0.1074922985149257
$ws.Cells.Item(7, "T").Value = 0.1074922985149257

# Row 8
$ws.Cells.Item(8, "A").Value = "sCs"
$ws.Cells.Item(8, "B").Value = "Egf"
$ws.Cells.Item(8, "C").Value = "Egfr"
$ws.Cells.Item(8, "D").Value = "ECs"
$ws.Cells.Item(8, "E").Value = 3
$ws.Cells.Item(8, "F").Value = 1
$ws.Cells.Item(8, "G").Value = 0.3380806666666666
$ws.Cells.Item(8, "H").Value = 1.014242
$ws.Cells.Item(8, "I").Value = 0.4234616956721602
$ws.Cells.Item(8, "J").Value = 0.4234616956721601
$ws.Cells.Item(8, "K").Value = 3
$ws.Cells.Item(8, "L").Value = 1
$ws.Cells.Item(8, "M").Value = 1.307106666666667
$ws.Cells.Item(8, "N").Value = 3.92132
$ws.Cells.Item(8, "O").Value = 0.01256263154946851
$ws.Cells.Item(8, "P").Value = 0.01256263154946851
$ws.Cells.Item(8, "Q").Value = 0.4419074932711111
$ws.Cells.Item(8, "R").Value = 3.97716743944
$ws.Cells.Item(8, "S").Value = 0.005319793258042512
$ws.Cells.Item(8, "T").Value = 0.005319793258042512

# Row 9
$ws.Cells.Item(9, "A").Value = "sCs"
$ws.Cells.Item(9, "B").Value = "Egf"
$ws.Cells.Item(9, "C").Value = "Egfr"
$ws.Cells.Item(9, "D").Value = "FAPs"
$ws.Cells.Item(9, "E").Value = 3
$ws.Cells.Item(9, "F").Value = 1
$ws.Cells.Item(9, "G").Value = 0.3380806666666666
$ws.Cells.Item(9, "H").Value = 1.014242
$ws.Cells.Item(9, "I").Value = 0.4234616956721602
$ws.Cells.Item(9, "J").Value = 0.4234616956721601
$ws.Cells.Item(9, "K").Value = 3
$ws.Cells.Item(9, "L").Value = 1
$ws.Cells.Item(9, "M").Value = 80.22623699999998
$ws.Cells.Item(9, "N").Value = 240.678711
$ws.Cells.Item(9, "O").Value = 0.77105616682495
$ws.Cells.Item(9, "P").Value = 0.77105616682495
$ws.Cells.Item(9, "Q").Value = 27.12293968911799
$ws.Cells.Item(9, "R").Value = 244.1064572020619
$ws.Cells.Item(9, "S").Value = 0.3265127518621694
$ws.Cells.Item(9, "T").Value = 0.3265127518621693

# Row 10
$ws.Cells.Item(10, "A").Value = "sCs"
$ws.Cells.Item(10, "B").Value = "Egf"
$ws.Cells.Item(10, "C").Value = "Egfr"
$ws.Cells.Item(10, "D").Value = "sCs"
$ws.Cells.Item(10, "E").Value = 3
$ws.Cells.Item(10, "F").Value = 1
$ws.Cells.Item(10, "G").Value = 0.3380806666666666
$ws.Cells.Item(10, "H").Value = 1.014242
$ws.Cells.Item(10, "I").Value = 0.4234616956721602
$ws.Cells.Item(10, "J").Value = 0.4234616956721601
$ws.Cells.Item(10, "K").Value = 3
$ws.Cells.Item(10, "L").Value = 1
$ws.Cells.Item(10, "M").Value = 22.51385866666667
$ws.Cells.Item(10, "N").Value = 67.541576
$ws.Cells.Item(10, "O").Value = 0.2163812016255815
$ws.Cells.Item(10, "P").Value = 0.2163812016255815
$ws.Cells.Item(10, "Q").Value = 7.611500347265777
$ws.Cells.Item(10, "R").Value = 68.503503125392
$ws.Cells.Item(10, "S").Value = 0.09162915055194831
$ws.Cells.Item(10, "T").Value = 0.0916291505519483
